# Auto-generated: apply updated market-price values per the commit diff.
# Values come from a scheduled market-data refresh; cells hold plain numbers (no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1415.3438
$ws.Range("J17").Value = 1104.129
$ws.Range("L17").Value = 3312.387
$ws.Range("N17").Value = -3648.387

$ws.Range("H38").Value = 1486.3334
$ws.Range("I38").Value = 1257.75
$ws.Range("J38").Value = 1943.5
$ws.Range("K38").Value = 3773.25
$ws.Range("L38").Value = 5830.5
$ws.Range("M38").Value = -3401.25
$ws.Range("N38").Value = -6574.5

$ws.Range("H82").Value = 4339.8
$ws.Range("I82").Value = 4339.8
$ws.Range("K82").Value = 13019.4
$ws.Range("M82").Value = -12613.4

$ws.Range("H85").Value = 4339.8
$ws.Range("I85").Value = 4339.8
$ws.Range("K85").Value = 13019.4
$ws.Range("M85").Value = -11615.4

$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120

$ws.Range("H138").Value = 3219.6
$ws.Range("I138").Value = 3819.0625
$ws.Range("J138").Value = 2153.889
$ws.Range("K138").Value = 11457.1875
$ws.Range("L138").Value = 6461.667
$ws.Range("M138").Value = -6317.1875
$ws.Range("N138").Value = -16741.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2877.879
$ws.Range("I32").Value = 1819.7693
$ws.Range("K32").Value = 1819.7693
$ws.Range("M32").Value = -1532.7693

$ws.Range("H74").Value = 3008.9473
$ws.Range("I74").Value = 2973.0625
$ws.Range("K74").Value = 2973.0625
$ws.Range("M74").Value = -2099.0625

$ws.Range("H77").Value = 3008.9473
$ws.Range("I77").Value = 2973.0625
$ws.Range("K77").Value = 14865.3125
$ws.Range("M77").Value = -10497.3125

$ws.Range("H110").Value = 1659.8
$ws.Range("I110").Value = 1659.8
$ws.Range("K110").Value = 1659.8
$ws.Range("M110").Value = 385.2

$ws.Range("H122").Value = 1815.5834
$ws.Range("I122").Value = 1845.2778
$ws.Range("K122").Value = 5535.8334
$ws.Range("M122").Value = -3085.8334

$ws.Range("H132").Value = 1733.3182
$ws.Range("I132").Value = 909.6429000000001
$ws.Range("K132").Value = 2728.9287
$ws.Range("M132").Value = -198.9287000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1958.5714
$ws.Range("I107").Value = 1587.8182
$ws.Range("K107").Value = 1587.8182
$ws.Range("M107").Value = 332.1818000000001

$ws.Range("H134").Value = 8054.607
$ws.Range("I134").Value = 8549.333000000001
$ws.Range("K134").Value = 25647.999
$ws.Range("M134").Value = -23112.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 932.6667
$ws.Range("I16").Value = 932.6667
$ws.Range("K16").Value = 932.6667
$ws.Range("M16").Value = -645.6667

$ws.Range("H31").Value = 3793.9333
$ws.Range("I31").Value = 1101.2858
$ws.Range("J31").Value = 6150
$ws.Range("K31").Value = 1101.2858
$ws.Range("L31").Value = 6150
$ws.Range("M31").Value = -806.2858000000001
$ws.Range("N31").Value = -6740

$ws.Range("H34").Value = 3793.9333
$ws.Range("I34").Value = 1101.2858
$ws.Range("J34").Value = 6150
$ws.Range("K34").Value = 1101.2858
$ws.Range("L34").Value = 6150
$ws.Range("M34").Value = -899.2858000000001
$ws.Range("N34").Value = -6554

$ws.Range("H41").Value = 12300.917
$ws.Range("J41").Value = 29000
$ws.Range("L41").Value = 29000
$ws.Range("N41").Value = -29856

$ws.Range("H86").Value = 1095.375
$ws.Range("I86").Value = 1109.1428
$ws.Range("K86").Value = 1109.1428
$ws.Range("M86").Value = 13.85719999999992

$ws.Range("H89").Value = 1095.375
$ws.Range("I89").Value = 1109.1428
$ws.Range("K89").Value = 5545.714
$ws.Range("M89").Value = 70.28600000000006

$ws.Range("H99").Value = 2284.375
$ws.Range("I99").Value = 2155
$ws.Range("K99").Value = 2155
$ws.Range("M99").Value = -657

$ws.Range("H113").Value = 932.6667
$ws.Range("I113").Value = 932.6667
$ws.Range("K113").Value = 932.6667
$ws.Range("M113").Value = 1237.3333

$ws.Range("H126").Value = 2284.375
$ws.Range("I126").Value = 2155
$ws.Range("K126").Value = 6465
$ws.Range("M126").Value = -3995

$ws.Range("H134").Value = 3042.5
$ws.Range("I134").Value = 2561.1
$ws.Range("K134").Value = 7683.299999999999
$ws.Range("M134").Value = -5148.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11710.429
$ws.Range("I87").Value = 657.6667
$ws.Range("K87").Value = 1973.0001
$ws.Range("M87").Value = -725.0001

$ws.Range("H90").Value = 11710.429
$ws.Range("I90").Value = 657.6667
$ws.Range("K90").Value = 5919.0003
$ws.Range("M90").Value = 320.9997000000003

$ws.Range("H131").Value = 9448200
$ws.Range("I131").Value = 62500588
$ws.Range("J131").Value = 16664.71
$ws.Range("K131").Value = 187501764
$ws.Range("L131").Value = 49994.13
$ws.Range("M131").Value = -187496724
$ws.Range("N131").Value = -60074.13

$ws.Range("H139").Value = 7562.8125
$ws.Range("I139").Value = 7940.3335
$ws.Range("J139").Value = 1900
$ws.Range("K139").Value = 23821.0005
$ws.Range("L139").Value = 5700
$ws.Range("M139").Value = -18681.0005
$ws.Range("N139").Value = -15980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2516
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 2595.5557
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 2595.5557
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -4591.5557

$ws.Range("H83").Value = 2516
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 2595.5557
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 12977.7785
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -22961.7785

$ws.Range("H102").Value = 2362.9473
$ws.Range("I102").Value = 2840.25
$ws.Range("K102").Value = 2840.25
$ws.Range("M102").Value = -1218.25

$ws.Range("H126").Value = 74716.14
$ws.Range("I126").Value = 3637.4546
$ws.Range("K126").Value = 10912.3638
$ws.Range("M126").Value = -8442.363799999999

$ws.Range("H132").Value = 4205
$ws.Range("I132").Value = 3643.2
$ws.Range("K132").Value = 10929.6
$ws.Range("M132").Value = -8399.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1117.25
$ws.Range("I22").Value = 480
$ws.Range("J22").Value = 1329.6666
$ws.Range("K22").Value = 480
$ws.Range("L22").Value = 1329.6666
$ws.Range("M22").Value = -185
$ws.Range("N22").Value = -1919.6666

$ws.Range("H27").Value = 1117.25
$ws.Range("I27").Value = 480
$ws.Range("J27").Value = 1329.6666
$ws.Range("K27").Value = 480
$ws.Range("L27").Value = 1329.6666
$ws.Range("M27").Value = -373
$ws.Range("N27").Value = -1543.6666

$ws.Range("H40").Value = 5022.407
$ws.Range("I40").Value = 3039.8
$ws.Range("K40").Value = 3039.8
$ws.Range("M40").Value = -2903.8

$ws.Range("H122").Value = 5913.9546
$ws.Range("I122").Value = 4283.8335
$ws.Range("K122").Value = 12851.5005
$ws.Range("M122").Value = -10401.5005

$ws.Range("H132").Value = 1771.8235
$ws.Range("I132").Value = 1044.5714
$ws.Range("K132").Value = 3133.7142
$ws.Range("M132").Value = -603.7142000000003

$ws.Range("H136").Value = 3300.0293
$ws.Range("I136").Value = 2148.0435
$ws.Range("K136").Value = 6444.130500000001
$ws.Range("M136").Value = -3894.130500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H107").Value = 921.9
$ws.Range("I107").Value = 748.5714
$ws.Range("J107").Value = 1326.3334
$ws.Range("K107").Value = 2245.7142
$ws.Range("L107").Value = 3979.0002
$ws.Range("M107").Value = -325.7142000000003
$ws.Range("N107").Value = -7819.0002

$ws.Range("H122").Value = 42911.42
$ws.Range("I122").Value = 78952
$ws.Range("K122").Value = 236856
$ws.Range("M122").Value = -234406

$ws.Range("H126").Value = 4928.7085
$ws.Range("I126").Value = 4405.375
$ws.Range("J126").Value = 5975.375
$ws.Range("K126").Value = 13216.125
$ws.Range("L126").Value = 17926.125
$ws.Range("M126").Value = -10746.125
$ws.Range("N126").Value = -22866.125

$ws.Range("H132").Value = 1900.2941
$ws.Range("I132").Value = 1108.2307
$ws.Range("J132").Value = 4474.5
$ws.Range("K132").Value = 3324.6921
$ws.Range("L132").Value = 13423.5
$ws.Range("M132").Value = -794.6921000000002
$ws.Range("N132").Value = -18483.5
